$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) and E (Volume) to Text format so that
# numeric-looking strings (e.g. "11.67", "102.90", "42.502.40") are
# preserved exactly as text, matching the source inline-string cells.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '42.502.40'
$ws.Range("E2").Value = '  -2.63%  '

# Row 3
$ws.Range("D3").Value = '2.229.94'
$ws.Range("E3").Value = '  -1.80%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").Value = '112.59'
$ws.Range("E5").Value = '  -6.78%  '

# Row 6
$ws.Range("D6").Value = '297.85'
$ws.Range("E6").Value = '  +11.92%  '

# Row 7
$ws.Range("D7").Value = '0.620'
$ws.Range("E7").Value = '  -3.82%  '

# Row 8
$ws.Range("E8").Value = '  -0.37%  '

# Row 9
$ws.Range("D9").Value = '0.605'
$ws.Range("E9").Value = '  -2.94%  '

# Row 10
$ws.Range("D10").Value = '44.35'
$ws.Range("E10").Value = '  -7.26%  '

# Row 11
$ws.Range("D11").Value = '0.0913'
$ws.Range("E11").Value = '  -3.27%  '

# Row 12
$ws.Range("D12").Value = '54.55'
$ws.Range("E12").Value = '  +0.90%  '

# Row 13
$ws.Range("D13").Value = '8.80'
$ws.Range("E13").Value = '  -6.71%  '

# Row 14
$ws.Range("D14").Value = '1.02'
$ws.Range("E14").Value = '  +13.60%  '

# Row 15
$ws.Range("D15").Value = '0.104'
$ws.Range("E15").Value = '  -1.99%  '

# Row 16
$ws.Range("D16").Value = '15.09'
$ws.Range("E16").Value = '  -3.22%  '

# Row 17
$ws.Range("D17").Value = '2.563.35'
$ws.Range("E17").Value = '  -1.93%  '

# Row 18
$ws.Range("D18").Value = '2.245.02'
$ws.Range("E18").Value = '  -1.15%  '

# Row 19
$ws.Range("D19").Value = '42.502.63'
$ws.Range("E19").Value = '  -2.60%  '

# Row 20
$ws.Range("D20").Value = '7.26'
$ws.Range("E20").Value = '  +4.91%  '

# Row 21
$ws.Range("E21").Value = '  -3.72%  '

# Row 22
$ws.Range("D22").Value = '73.62'
$ws.Range("E22").Value = '  +2.02%  '

# Row 23
$ws.Range("D23").Value = '3.55'
$ws.Range("E23").Value = '  +22.29%  '

# Row 24
$ws.Range("D24").Value = '2.39'
$ws.Range("E24").Value = '  -1.73%  '

# Row 25
$ws.Range("D25").Value = '230.41'
$ws.Range("E25").Value = '  -1.75%  '

# Row 26
$ws.Range("D26").Value = '9.26'
$ws.Range("E26").Value = '  -2.77%  '

# Row 27
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '11.67'
$ws.Range("E27").Value = '  -1.47%  '

# Row 28
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -1.55%  '

# Row 29
$ws.Range("D29").Value = '38.95'
$ws.Range("E29").Value = '  -8.88%  '

# Row 30
$ws.Range("E30").Value = '  -0.48%  '

# Row 31
$ws.Range("E31").Value = '  -4.26%  '

# Row 32
$ws.Range("D32").Value = '173.99'
$ws.Range("E32").Value = '  -0.16%  '

# Row 33
$ws.Range("E33").Value = '  -2.34%  '

# Row 34
$ws.Range("D34").Value = '0.0895'
$ws.Range("E34").Value = '  -3.23%  '

# Row 35
$ws.Range("D35").Value = '5.70'
$ws.Range("E35").Value = '  -0.96%  '

# Row 36
$ws.Range("D36").Value = '5.16'
$ws.Range("E36").Value = '  +11.58%  '

# Row 37
$ws.Range("D37").Value = '4.39'
$ws.Range("E37").Value = '  +2.68%  '

# Row 38
$ws.Range("E38").Value = '  -3.47%  '

# Row 39
$ws.Range("E39").Value = '  -1.92%  '

# Row 40
$ws.Range("D40").Value = '0.105'
$ws.Range("E40").Value = '  -4.57%  '

# Row 41
$ws.Range("B41").Value = 'MultiversX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D41").Value = '74.39'
$ws.Range("E41").Value = '  +2.58%  '

# Row 42
$ws.Range("B42").Value = 'LidoDAOToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D42").Value = '2.42'
$ws.Range("E42").Value = '  -5.19%  '

# Row 43
$ws.Range("D43").Value = '0.235'
$ws.Range("E43").Value = '  -1.60%  '

# Row 44
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.04%  '

# Row 45
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").Value = '12.62'
$ws.Range("E45").Value = '  -8.39%  '

# Row 46
$ws.Range("E46").Value = '  -4.76%  '

# Row 47
$ws.Range("E47").Value = '  -5.72%  '

# Row 48
$ws.Range("D48").Value = '1.31'
$ws.Range("E48").Value = '  +2.73%  '

# Row 49
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = '1.68'
$ws.Range("E49").Value = '  +8.31%  '

# Row 50
$ws.Range("D50").Value = '8.59'
$ws.Range("E50").Value = '  +1.70%  '

# Row 51
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '102.90'
$ws.Range("E51").Value = '  -0.39%  '
